$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that needs to move
# from 45182 (2023-09-13) to 45184 (2023-09-15) for every data row (2-260).
$ws.Range("C2:C260").Value = 45184
